$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A82").Value = "2025-04-29 12:42:17"
$ws.Range("B82").Value = 261
